# Update column G ("K") values on Sheet1 with newly regenerated figures.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 4
    7  = 2
    8  = 4
    9  = 1
    10 = 4
    11 = 5
    12 = 2
    13 = 1
    14 = 1
    15 = 4
    16 = 1
    17 = 3
    18 = 7
    19 = 1
    20 = 4
    21 = 3
    22 = 6
    23 = 3
    24 = 5
    25 = 2
    26 = 2
    27 = 5
    28 = 6
    29 = 4
    30 = 4
    31 = 4
    32 = 6
    33 = 6
    34 = 1
    35 = 2
    36 = 3
}

foreach ($r in $newK.Keys) {
    $ws.Range("G$r").Value = $newK[$r]
}
